$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the currency code column (Cod. Moneda) from "US$" to "USD"
# for every data row in the table (rows 2-18, column F).
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Text -eq "US$") {
        $cell.Value = "USD"
    }
}
